$d = $word.ActiveDocument

$old = ": Daty kampanii używające Gwiazdozbiór Herkulesa 2022: 13-22 czerwca, 12-21 lipca, 10-19 sierpnia"
$new = "2022: Daty kampanii używające Gwiazdozbiór Herkulesa: 13-22 czerwca, 12-21 lipca, 10-19 sierpnia"

$d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
